# Registros_11-2025.xlsx — remove Taissa's 04/11 entry and fold the old
# 05/11 entry up (now missing its saída/intervalo punches), updating the
# dependent totals, dimensions and filter ranges to match.

$wb = $excel.ActiveWorkbook

# --- "Taissa" sheet -------------------------------------------------------
$wsTaissa = $wb.Worksheets.Item("Taissa")

# Drop row 2 (03/11 Seg) and the original row 3 (04/11 Ter) entirely; after
# each delete the rows below shift up, so deleting row 2 twice removes both.
$wsTaissa.Rows.Item(2).Delete()
$wsTaissa.Rows.Item(2).Delete()

# What used to be row 4 (05/11 Qua) is now row 2. Its punches are reduced to
# a single (later) "Entrada" time, with the rest of the day left blank and
# the computed hours zeroed out.
$wsTaissa.Range("C2").Value = "23:42"
$wsTaissa.Range("D2").Value = ""
$wsTaissa.Range("E2").Value = ""
$wsTaissa.Range("F2").Value = ""
$wsTaissa.Range("G2").Value = 0

# The TOTAL row (now row 3) reflects the same zeroed-out hours.
$wsTaissa.Range("G3").Value = 0

# Re-anchor the autofilter to the new (now 3-row) extent.
$wsTaissa.AutoFilterMode = $false
$wsTaissa.Range("A1:G3").AutoFilter()

# The hidden _FilterDatabase defined name for "Taissa" still points at the
# old $A$1:$G$5 extent; point it at the new $A$1:$G$3 extent.
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name() -eq "Taissa!_FilterDatabase") {
    $n.RefersTo = "='Taissa'!`$A`$1:`$G`$3"
  }
}

# --- "Resumo" sheet --------------------------------------------------------
# Taissa's total hours (summarised on the Resumo sheet) drop to zero too.
$wsResumo = $wb.Worksheets.Item("Resumo")
$wsResumo.Range("C3").Value = 0
